# Commit: "adicionado utils, pathlib e status atrasado para não enviados"
#
# For every data row whose "Situação" (column E) is still blank (i.e. the
# submission has neither been sent nor flagged any other way), mark it as
# "Atrasado" (late) with a white-on-orange highlight, matching the style
# already used for the other status labels ("Enviado" / "Envio Duplicado"
# use a white-on-dark-green highlight).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OLE (BGR) colors: Excel's Color properties take 0xBBGGRR, not 0xRRGGBB.
$white  = 16777215   # RGB(255,255,255) -> FFFFFF
$orange = 25855       # RGB(255,100,0)   -> FF6400

# Rows where column E (Situação) is currently empty.
$rows = @(2,3,4,5,6,8,9,11,12,13,14,15,16,18,19,22,23,24,25,26,27,29,31,32,33,35,37,38,39,40,41,42,43,44,45,47,48,49,50,51,53,54,55,56,58,60,61,62,63,67,68)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $cell.Value = "Atrasado"
    $cell.Font.Color = $white
    $cell.Interior.Color = $orange
    $cell.Interior.PatternColor = $orange
}
